$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header labels: _old -> _FV2410, _new -> _FV2504
$oldCols = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($oldCols[$i])_FV2410"
}

for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($oldCols[$i])_FV2504"
}

# Create table over the data range
$rng = $ws.Range("A1:U94")
$tbl = $ws.ListObjects.Add(1, $rng, 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
